$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "week 2" design tasks appended to the timeline table (rows 18-23).
$tasks = @(
    @{ Row = 18; Name = "Design the Wireframes(Lhub_Wireframe)"; Assigned = "Azmy" },
    @{ Row = 19; Name = "Design Use Case Diagrams(Lhub_usecase)"; Assigned = "Wafaa" },
    @{ Row = 20; Name = "Design the Entity Relationship Diagram (ERD)(Lhub_ERD)"; Assigned = "Aya Mohamed" },
    @{ Row = 21; Name = "Design the Class Diagram"; Assigned = "Farah" },
    @{ Row = 22; Name = "Design Data Flow Diagram(Lhub_DFD)"; Assigned = "Nada" },
    @{ Row = 23; Name = "Design Peer Review Sheet"; Assigned = "Aya Mohamed" }
)

foreach ($task in $tasks) {
    $r = $task.Row
    $ws.Range("A$r").Value = $task.Name
    $ws.Range("B$r").Value = $task.Assigned
    $ws.Range("C$r").Value = 43589
    $ws.Range("D$r").Value = 43596
    $ws.Range("E$r").Formula = "=(D$r-C$r)"
    $ws.Range("F$r").Value = 100
    $ws.Range("G$r").Value = "week 2"
}
